$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Province names: abbreviation -> full name
$ws.Range("C2").Value = "Alberta"
$ws.Range("C3").Value = "Colombie-Britannique"

# Column C (PROVINCE) widened to fit the longer text
$ws.Columns("C").ColumnWidth = 17.8

# Selection moved to D7
$ws.Range("D7").Select()
